$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 647.1429
$ws.Range("I4").Value = 705
$ws.Range("K4").Value = 705
$ws.Range("M4").Value = -591

$ws.Range("H21").Value = 13412.8
$ws.Range("I21").Value = 9908.211
$ws.Range("J21").Value = 80000
$ws.Range("K21").Value = 9908.211
$ws.Range("L21").Value = 80000
$ws.Range("M21").Value = -9440.211
$ws.Range("N21").Value = -80936

$ws.Range("H23").Value = 13412.8
$ws.Range("I23").Value = 9908.211
$ws.Range("J23").Value = 80000
$ws.Range("K23").Value = 9908.211
$ws.Range("L23").Value = 80000
$ws.Range("M23").Value = -9674.211
$ws.Range("N23").Value = -80468

$ws.Range("H29").Value = 638.6
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws.Range("H33").Value = 170.96297
$ws.Range("I33").Value = 145.47368
$ws.Range("K33").Value = 145.47368
$ws.Range("M33").Value = 83.52632

$ws.Range("H51").Value = 2992
$ws.Range("I51").Value = 2980
$ws.Range("J51").Value = 2994.4
$ws.Range("K51").Value = 2980
$ws.Range("L51").Value = 2994.4
$ws.Range("M51").Value = -2496
$ws.Range("N51").Value = -3962.4

$ws.Range("H125").Value = 1438.1666
$ws.Range("I125").Value = 900
$ws.Range("J125").Value = 1976.3334
$ws.Range("K125").Value = 8100
$ws.Range("L125").Value = 17787.0006
$ws.Range("M125").Value = -5640
$ws.Range("N125").Value = -22707.0006

$ws.Range("H135").Value = 2837.25
$ws.Range("I135").Value = 1041.3334
$ws.Range("J135").Value = 8225
$ws.Range("K135").Value = 9372.0006
$ws.Range("L135").Value = 74025
$ws.Range("M135").Value = -6837.000599999999
$ws.Range("N135").Value = -79095


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 928.6667
$ws.Range("I45").Value = 866.3333
$ws.Range("J45").Value = 1053.3334
$ws.Range("K45").Value = 866.3333
$ws.Range("L45").Value = 1053.3334
$ws.Range("M45").Value = -489.3333
$ws.Range("N45").Value = -1807.3334

$ws.Range("H63").Value = 5499.6665
$ws.Range("I63").Value = 2500
$ws.Range("K63").Value = 2500
$ws.Range("M63").Value = -1814

$ws.Range("H66").Value = 5499.6665
$ws.Range("I66").Value = 2500
$ws.Range("K66").Value = 12500
$ws.Range("M66").Value = -9068

$ws.Range("H122").Value = 2494.3333
$ws.Range("J122").Value = 2491
$ws.Range("L122").Value = 7473
$ws.Range("N122").Value = -12373

$ws.Range("H132").Value = 3802.125
$ws.Range("I132").Value = 3893.6
$ws.Range("J132").Value = 3649.6667
$ws.Range("K132").Value = 11680.8
$ws.Range("L132").Value = 10949.0001
$ws.Range("M132").Value = -9150.8
$ws.Range("N132").Value = -16009.0001


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 1916.4482
$ws.Range("I86").Value = 1399.5834
$ws.Range("J86").Value = 2281.2942
$ws.Range("K86").Value = 1399.5834
$ws.Range("L86").Value = 2281.2942
$ws.Range("M86").Value = -276.5834
$ws.Range("N86").Value = -4527.2942

$ws.Range("H89").Value = 1916.4482
$ws.Range("I89").Value = 1399.5834
$ws.Range("J89").Value = 2281.2942
$ws.Range("K89").Value = 6997.916999999999
$ws.Range("L89").Value = 11406.471
$ws.Range("M89").Value = -1381.916999999999
$ws.Range("N89").Value = -22638.471

$ws.Range("H99").Value = 1284
$ws.Range("I99").Value = 1252.8
$ws.Range("K99").Value = 1252.8
$ws.Range("M99").Value = 245.2

$ws.Range("H126").Value = 1284
$ws.Range("I126").Value = 1252.8
$ws.Range("K126").Value = 3758.4
$ws.Range("M126").Value = -1288.4


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1173
$ws.Range("I4").Value = 900.5
$ws.Range("J4").Value = 1282
$ws.Range("K4").Value = 2701.5
$ws.Range("L4").Value = 3846
$ws.Range("M4").Value = -2589.5
$ws.Range("N4").Value = -4070

$ws.Range("H14").Value = 549.8333
$ws.Range("I14").Value = 549.8333
$ws.Range("K14").Value = 1649.4999
$ws.Range("M14").Value = -1476.4999

$ws.Range("H50").Value = 11488.889
$ws.Range("I50").Value = 16916.666
$ws.Range("J50").Value = 633.3333
$ws.Range("K50").Value = 50749.99800000001
$ws.Range("L50").Value = 1899.9999
$ws.Range("M50").Value = -50268.99800000001
$ws.Range("N50").Value = -2861.9999

$ws.Range("H53").Value = 11488.889
$ws.Range("I53").Value = 16916.666
$ws.Range("J53").Value = 633.3333
$ws.Range("K53").Value = 50749.99800000001
$ws.Range("L53").Value = 1899.9999
$ws.Range("M53").Value = -50268.99800000001
$ws.Range("N53").Value = -2861.9999

$ws.Range("H98").Value = 1698.1818
$ws.Range("I98").Value = 998
$ws.Range("J98").Value = 2538.4
$ws.Range("K98").Value = 2994
$ws.Range("L98").Value = 7615.200000000001
$ws.Range("M98").Value = -1496
$ws.Range("N98").Value = -10611.2

$ws.Range("H106").Value = 5285.7144
$ws.Range("J106").Value = 5285.7144
$ws.Range("L106").Value = 15857.1432
$ws.Range("N106").Value = -17749.1432


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 20125
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 22857.143
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 22857.143
$ws.Range("M43").Value = -849
$ws.Range("N43").Value = -23159.143

$ws.Range("H46").Value = 23500
$ws.Range("J46").Value = 23500
$ws.Range("L46").Value = 23500
$ws.Range("N46").Value = -23812

$ws.Range("H57").Value = 15040
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 15040
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 15040
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -16680

$ws.Range("H80").Value = 173129.42
$ws.Range("J80").Value = 241041.2
$ws.Range("L80").Value = 241041.2
$ws.Range("N80").Value = -243037.2

$ws.Range("H83").Value = 173129.42
$ws.Range("J83").Value = 241041.2
$ws.Range("L83").Value = 1205206
$ws.Range("N83").Value = -1215190

$ws.Range("H102").Value = 1158.8182
$ws.Range("I102").Value = 851.7647
$ws.Range("K102").Value = 851.7647
$ws.Range("M102").Value = 770.2353

$ws.Range("H113").Value = 2904.7144
$ws.Range("I113").Value = 4142.3335
$ws.Range("J113").Value = 1976.5
$ws.Range("K113").Value = 4142.3335
$ws.Range("L113").Value = 1976.5
$ws.Range("M113").Value = -1972.3335
$ws.Range("N113").Value = -6316.5

$ws.Range("H122").Value = 5004.636
$ws.Range("I122").Value = 6336.143
$ws.Range("J122").Value = 2674.5
$ws.Range("K122").Value = 19008.429
$ws.Range("L122").Value = 8023.5
$ws.Range("M122").Value = -16558.429
$ws.Range("N122").Value = -12923.5

$ws.Range("H126").Value = 3816.5
$ws.Range("I126").Value = 3849.75
$ws.Range("K126").Value = 11549.25
$ws.Range("M126").Value = -9079.25


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws.Range("H46").Value = 1011.1111
$ws.Range("I46").Value = 1300
$ws.Range("J46").Value = 780
$ws.Range("K46").Value = 1300
$ws.Range("L46").Value = 780
$ws.Range("M46").Value = -1112
$ws.Range("N46").Value = -1156

$ws.Range("H61").Value = 660.4211
$ws.Range("I61").Value = 583
$ws.Range("J61").Value = 877.2
$ws.Range("K61").Value = 583
$ws.Range("L61").Value = 877.2
$ws.Range("M61").Value = -381
$ws.Range("N61").Value = -1281.2

$ws.Range("H113").Value = 660.4211
$ws.Range("I113").Value = 583
$ws.Range("J113").Value = 877.2
$ws.Range("K113").Value = 583
$ws.Range("L113").Value = 877.2
$ws.Range("M113").Value = 1587
$ws.Range("N113").Value = -5217.2

$ws.Range("H132").Value = 8142.619
$ws.Range("I132").Value = 13044.728
$ws.Range("K132").Value = 39134.18399999999
$ws.Range("M132").Value = -36604.18399999999

